$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Row 5 (ID=4): STATO flips from TODO -> DONE.
#    Copy the "Good" (DONE) cell format from E2 so the existing style entry
#    (green fill / center+wrap) is reused instead of creating a new one.
# ---------------------------------------------------------------------------
$ws.Range("E5").Value = "DONE"
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Row 7 (ID=6): TIPO Bug -> Task, DESCRIZIONE wording tweak, STATO
#    TODO -> DONE, and the row shrinks because the new text wraps less.
# ---------------------------------------------------------------------------
$ws.Range("C7").Value = "Task"
$ws.Range("D7").Value = "La selezione dei settori non funziona ancora (passando da un SG all'altro, dovrebbe colorare di rosso quelli del SG corrente e, di bianco quelli relativi agli altri SG)"
$ws.Range("E7").Value = "DONE"
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Rows(7).RowHeight = 58

# ---------------------------------------------------------------------------
# 3) Row 8: was an empty placeholder row, becomes a new data row (ID=7).
#    Columns A-D/F already carry the correct placeholder formatting, so only
#    their values need to be populated. G8 needs both value and format
#    (copied from G2) since it had no cell at all before. E8 gets the new
#    "IN PROGRESS" status with a Neutral (amber) look.
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "SG Editor"
$ws.Range("C8").Value = "Bug"
$ws.Range("D8").Value = "Passando da una room all'altra, se ruoto o faccio zoom su una sembra che anche l'altra (o le altre) tengono in memoria qualche trasformazione, invece dovrebbero essere indipendneti."

$ws.Range("E8").Value = "IN PROGRESS"
$ws.Range("E8").Style = "Neutral"
$ws.Range("E8").HorizontalAlignment = -4108   # xlCenter
$ws.Range("E8").VerticalAlignment = -4108     # xlCenter
$ws.Range("E8").WrapText = $true

$ws.Range("G8").Value = "Alta"
$ws.Range("G2").Copy() | Out-Null
$ws.Range("G8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Rows(8).RowHeight = 72.5

# ---------------------------------------------------------------------------
# 4) Column E widened slightly (best-fit for "IN PROGRESS").
# ---------------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 10.9

# ---------------------------------------------------------------------------
# 5) View state: frozen pane / active selection follow the new data extent.
# ---------------------------------------------------------------------------
$ws.Range("H7").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.SplitRow = 1
$excel.ActiveWindow.SplitColumn = 7
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("G8").Select()

$excel.CutCopyMode = $false
